$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F5 was "quantity" -> becomes "product_id"
$ws.Range("F5").Value = "product_id"
# F6 is new -> "quantity"
$ws.Range("F6").Value = "quantity"
# F7 is new -> "status" (new shared string)
$ws.Range("F7").Value = "status"

# Final selection ends on F7
$ws.Range("F7").Select()
